# AFDP-7308 Combine Transcribe and OCR processing into a single media processing module
#
# This updates the Drools "Transcribe Rules" decision table so that it
# references the new combined MediaEngine model instead of the old
# Transcribe-only model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Import statement: switch the imported class from the Transcribe model to
# the new MediaEngine model.
$ws.Range("D4").Value = "com.armedia.acm.services.mediaengine.model.MediaEngine"

# RuleTable header binding: bind $transcribe to the MediaEngine type instead
# of the Transcribe type.
$ws.Range("C17").Value = "`$transcribe: MediaEngine"
